# Sprint Backlog update - "New Sprint has started!"
#
# Folha1 is a small Scrum task board (B6:E9 header + two task rows).
# Column C ("In progress") lists who is currently working each task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 1 ("1º functionality") now has more people helping out: Joana and
# Leticia joined Guilherme on this task.
$ws.Range("C8").Value = "Guilherme Joana         Leticia"

# Task 2 ("2º functionality") now has Samuel assigned to it.
$ws.Range("C9").Value = "Samuel"

# A new, still-empty row was started right below the table (row 10) and
# pre-formatted with an underlined font, ready for the next entry.
$ws.Range("C10").Font.Underline = $true

# Leave the selection on the newly-touched cell, same as the live editor.
[void]$ws.Range("C10").Select()
